$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-number-looking price values in column D to remain text
# (matches the source data which stores them as inline strings)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

# Apply updated values
$ws.Range("D2").Value = "42.398.93"
$ws.Range("E2").Value = "  +1.12%  "
$ws.Range("D3").Value = "2.290.90"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "322.46"
$ws.Range("E5").Value = "  +2.31%  "
$ws.Range("D6").Value = "104.22"
$ws.Range("E6").Value = "  +1.90%  "
$ws.Range("D7").Value = "0.628"
$ws.Range("E7").Value = "  +0.47%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +1.03%  "
$ws.Range("E10").Value = "  +3.17%  "
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("D12").Value = "8.56"
$ws.Range("E12").Value = "  +3.62%  "
$ws.Range("D13").Value = "0.106"
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("D14").Value = "0.967"
$ws.Range("E14").Value = "  +0.83%  "
$ws.Range("D15").Value = "15.26"
$ws.Range("E15").Value = "  +0.48%  "
$ws.Range("D16").Value = "2.639.13"
$ws.Range("D17").Value = "2.280.22"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").Value = "42.371.16"
$ws.Range("E18").Value = "  +1.50%  "
$ws.Range("D19").Value = "7.47"
$ws.Range("E19").Value = "  -0.77%  "
$ws.Range("E20").Value = "  +0.25%  "
$ws.Range("D21").Value = "13.28"
$ws.Range("E21").Value = "  +34.26%  "
$ws.Range("D22").Value = "73.28"
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("D23").Value = "3.60"
$ws.Range("E23").Value = "  +1.71%  "
$ws.Range("D24").Value = "269.19"
$ws.Range("E24").Value = "  -5.19%  "
$ws.Range("E25").Value = "  -1.66%  "
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("E27").Value = "  +1.40%  "
$ws.Range("D28").Value = "2.29"
$ws.Range("E28").Value = "  +0.20%  "
$ws.Range("D29").Value = "22.54"
$ws.Range("E29").Value = "  -1.68%  "
$ws.Range("D30").Value = "38.07"
$ws.Range("E30").Value = "  +10.34%  "
$ws.Range("D31").Value = "165.43"
$ws.Range("E31").Value = "  +1.62%  "
$ws.Range("D32").Value = "6.17"
$ws.Range("E32").Value = "  +6.46%  "
$ws.Range("D33").Value = "0.0882"
$ws.Range("E33").Value = "  +1.06%  "
$ws.Range("E34").Value = "  +0.58%  "
$ws.Range("E35").Value = "  -13.20%  "
$ws.Range("E36").Value = "  -1.27%  "
$ws.Range("E37").Value = "  +1.25%  "
$ws.Range("E38").Value = "  +2.35%  "
$ws.Range("D39").Value = "3.72"
$ws.Range("E39").Value = "  +3.71%  "
$ws.Range("D40").Value = "2.73"
$ws.Range("E40").Value = "  -5.45%  "
$ws.Range("E41").Value = "  +5.94%  "
$ws.Range("D42").Value = "70.00"
$ws.Range("E42").Value = "  +1.07%  "
$ws.Range("D43").Value = "96.31"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("E45").Value = "  +0.48%  "
$ws.Range("D46").Value = "12.41"
$ws.Range("E46").Value = "  +4.22%  "
$ws.Range("D47").Value = "81.41"
$ws.Range("E47").Value = "  +7.08%  "
$ws.Range("D48").Value = "113.11"
$ws.Range("E48").Value = "  -1.40%  "
$ws.Range("E49").Value = "  -0.88%  "
$ws.Range("E50").Value = "  -0.55%  "
$ws.Range("D51").Value = "1.587.83"
$ws.Range("E51").Value = "  +3.59%  "

# Remove the temporary text-number-format so styling matches the original (no style index on these cells)
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D7").ClearFormats()
$ws.Range("D12").ClearFormats()
$ws.Range("D13").ClearFormats()
$ws.Range("D14").ClearFormats()
$ws.Range("D15").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("D29").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("D48").ClearFormats()
